$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to underscored variants (keep Company at B1 unchanged)
$ws.Range("A1").Value = "Stock_Ticker"
$ws.Range("C1").Value = "Date_IPO"
$ws.Range("D1").Value = "Market_Capitalization"
$ws.Range("E1").Value = "Pet_Market_Segment"
$ws.Range("F1").Value = "Company_Size "

# Add a new row (row 7) for FRPT, copying formatting from the row above (row 6)
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

$ws.Range("A7").Value = "FRPT"
$ws.Range("B7").Value = "Freshpet, Inc."
$ws.Range("C7").Value = 41944
$ws.Range("D7").Value = "2.582B"
$ws.Range("E7").Value = "Pet Food"
$ws.Range("F7").Value = 789

# Update the active selection to the newly added cell, matching the saved view state
$ws.Range("F7").Select()
